$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.011440386912840304
$ws.Range("B1").Value = 0.49565056258015927
$ws.Range("C1").Value = 0.099958486403788641
$ws.Range("D1").Value = 295562842.63828558
$ws.Range("E1").Value = 30.007333053017643
$ws.Range("F1").Value = 54.884345453779154
$ws.Range("G1").Value = 0.37782256432557115
$ws.Range("H1").Value = 9.9999753171922663
$ws.Range("I1").Value = 0.015284193793147917
$ws.Range("J1").Value = 0.09347273948898363
$ws.Range("K1").Value = 0.97574148988616938
$ws.Range("L1").Value = 0.94357324530498876
$ws.Range("M1").Value = -0.48807107839730191
$ws.Range("A2").Value = 0.0027120401622155203
$ws.Range("B2").Value = -0.36126863009075266
$ws.Range("C2").Value = 0.048060698575486727
$ws.Range("D2").Value = 154280748.56200033
$ws.Range("E2").Value = 36.517931529860441
$ws.Range("F2").Value = 56.445421243458888
$ws.Range("G2").Value = 0.46565942347600586
$ws.Range("H2").Value = 0.63780196842218284
$ws.Range("I2").Value = 0.000034621105853416682
$ws.Range("J2").Value = 0.0000075097597061353193
$ws.Range("K2").Value = 0.92955544593543438
$ws.Range("L2").Value = 0.66523500227684074
$ws.Range("M2").Value = -2.0961005084669013
$ws.Range("A3").Value = 0.056819083470834864
$ws.Range("B3").Value = 13.302842815715305
$ws.Range("C3").Value = 0.099330058649682404
$ws.Range("D3").Value = 249407203.85844174
$ws.Range("E3").Value = 31.457138454781006
$ws.Range("F3").Value = 56.118050574322375
$ws.Range("G3").Value = 0.035981354888243465
$ws.Range("H3").Value = 9.0260406122547927
$ws.Range("I3").Value = 0.00366764977501033
$ws.Range("J3").Value = 0.098167824668004708
$ws.Range("K3").Value = 0.93027712136451557
$ws.Range("L3").Value = 0.29651402330969912
$ws.Range("M3").Value = 0.63637366557889907
$ws.Range("A4").Value = 0.02712689489805864
$ws.Range("B4").Value = 1.3714933383333792
$ws.Range("C4").Value = 0.087346457802709429
$ws.Range("D4").Value = 200605589.58113173
$ws.Range("E4").Value = 35.582107525914203
$ws.Range("F4").Value = 66.188666613795604
$ws.Range("G4").Value = 0.000046089882936936263
$ws.Range("H4").Value = 0.81719541243335991
$ws.Range("I4").Value = 0.10822344540187247
$ws.Range("J4").Value = 0.031460189547888893
$ws.Range("K4").Value = 0.99012870726610525
$ws.Range("L4").Value = 0.97107671097477977
$ws.Range("M4").Value = 0.99294200560831736
$ws.Range("A5").Value = 0.000029608641831456684
$ws.Range("B5").Value = -0.72281142087069294
$ws.Range("C5").Value = 0.058871789444177224
$ws.Range("D5").Value = 245534921.51240683
$ws.Range("E5").Value = 30.029572498396881
$ws.Range("F5").Value = 52.965500124559483
$ws.Range("G5").Value = 0.066205718916834716
$ws.Range("H5").Value = 8.1888278378098924
$ws.Range("I5").Value = 0.0023969090508505265
$ws.Range("J5").Value = 0.04086985569600337
$ws.Range("K5").Value = 0.86872676107072355
$ws.Range("L5").Value = -0.43521669109845296
$ws.Range("M5").Value = -0.61758975116199988
$ws.Range("A6").Value = -0.000013230881520724387
$ws.Range("B6").Value = -0.78489337783990254
$ws.Range("C6").Value = 0.08710743526393426
$ws.Range("D6").Value = 206166475.98038509
$ws.Range("E6").Value = 35.184319395031473
$ws.Range("F6").Value = 55.432560546964417
$ws.Range("G6").Value = 0.043194295876404346
$ws.Range("H6").Value = 0.22490550034772616
$ws.Range("I6").Value = 0.0047791216147433074
$ws.Range("J6").Value = 0.071737755819651625
$ws.Range("K6").Value = 0.97053154606191161
$ws.Range("L6").Value = -0.70522118356909402
$ws.Range("M6").Value = -0.53451576285609681
$ws.Range("A7").Value = 0.040390256405453226
$ws.Range("B7").Value = 0.3191779890306865
$ws.Range("C7").Value = 0.046889349217833505
$ws.Range("D7").Value = 177945685.1172958
$ws.Range("E7").Value = 38.68527965278934
$ws.Range("F7").Value = 54.389723143295861
$ws.Range("G7").Value = 0.038324829743957262
$ws.Range("H7").Value = 0.0075191169398032701
$ws.Range("I7").Value = 0.15609545590090443
$ws.Range("J7").Value = 0.0029227286724083251
$ws.Range("K7").Value = 0.99684532175112639
$ws.Range("L7").Value = 0.93826073361584261
$ws.Range("M7").Value = 0.79815707058722574
$ws.Range("A8").Value = 0.035433609589332747
$ws.Range("B8").Value = 0.2248120953199865
$ws.Range("C8").Value = 0.086814305824355842
$ws.Range("D8").Value = 313736434.43527311
$ws.Range("E8").Value = 34.293200552222373
$ws.Range("F8").Value = 62.606961317250672
$ws.Range("G8").Value = 0.13340528472058127
$ws.Range("H8").Value = 0.012200464342299471
$ws.Range("I8").Value = 0.028267980053705608
$ws.Range("J8").Value = 0.0039927307704885184
$ws.Range("K8").Value = 0.9974947501641046
$ws.Range("L8").Value = 0.98018026782951018
$ws.Range("M8").Value = 0.88487625956443938
$ws.Range("K9").Value = 0.063869368528193804
$ws.Range("A10").Value = 0.0075232950610628815
$ws.Range("B10").Value = 1.2141083972404731
$ws.Range("C10").Value = 0.027874915534674586
$ws.Range("D10").Value = 100000359.20397583
$ws.Range("E10").Value = 40.309714392356852
$ws.Range("F10").Value = 45.860540858324214
$ws.Range("G10").Value = 0.0015024779318060649
$ws.Range("H10").Value = 0.0000023837407607480294
$ws.Range("I10").Value = 0.8986371888176754
$ws.Range("J10").Value = 0.099882204667615909
$ws.Range("K10").Value = 0.99116143178843219
$ws.Range("L10").Value = 0.99279815658694825
$ws.Range("M10").Value = 0.96116099727156035
$ws.Range("A11").Value = 0.22762795461504831
$ws.Range("B11").Value = 4.2535265926058408
$ws.Range("C11").Value = 0.05583387189061783
$ws.Range("D11").Value = 192570108.45417479
$ws.Range("E11").Value = 30.000123938224252
$ws.Range("F11").Value = 59.457651373573249
$ws.Range("G11").Value = 2.2899629476418166
$ws.Range("H11").Value = 9.9973668218250804
$ws.Range("I11").Value = 0.053843188760794962
$ws.Range("J11").Value = 0.00052842847193436523
$ws.Range("K11").Value = 0.94492769577785629
$ws.Range("L11").Value = 0.95166691727919894
$ws.Range("M11").Value = 0.84653415764886653
$ws.Range("A12").Value = 0.11666514709222116
$ws.Range("B12").Value = 3.0782162999939455
$ws.Range("C12").Value = 0.099802561563821363
$ws.Range("D12").Value = 978225455.81800199
$ws.Range("E12").Value = 30.020264712271857
$ws.Range("F12").Value = 60.874625319891393
$ws.Range("G12").Value = 0.69823817773265762
$ws.Range("H12").Value = 5.2942953206266168
$ws.Range("I12").Value = 0.00062410739733497573
$ws.Range("J12").Value = 0.0077935959548451738
$ws.Range("K12").Value = 0.77939519750869957
$ws.Range("L12").Value = 0.77816846083135593
$ws.Range("M12").Value = 0.88166915610345598
$ws.Range("K13").Value = 0.10601660076756791
$ws.Range("K14").Value = 0.072690375408132502
$ws.Range("K15").Value = 0.054130836864196841
$ws.Range("K16").Value = 0.35146614768476769

Write-Output "Updated $($ws.Name) with smoothed values"
